$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-9 (Ano, total_item, Evolucao (%))
$data = @(
    @(2018, 700197.33, $null),
    @(2019, 1741717.08, 148.746604046605),
    @(2020, 3313857.03, 90.26379588583926),
    @(2021, 6564496.12, 98.09231540685992),
    @(2022, 7360567.04, 12.12691584316148),
    @(2023, 6932834.57, -5.811134762791315),
    @(2024, 9425497.51, 35.95445578330019),
    @(2025, 3265484.16, -65.35478199919443)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    if ($item[2] -ne $null) {
        $ws.Cells.Item($row, 3).Value = $item[2]
    } else {
        $ws.Cells.Item($row, 3).Value = $null
    }
    $row++
}
